$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Shift the existing "ConceptScheme table" (currently rows 7-11) down to rows 14-18
# by inserting 7 new rows above it (rows 4-10). Row formatting (bold labels in column A,
# wrap style in column C) is automatically extended into the newly inserted rows by Excel.
$ws.Rows("4:10").Insert()

# --- Update existing metadata row 3 (dct:description text changes) ---
$ws.Range("B3").Value = "Controlled vocabulary for the status of data sources in TANGENT."

# --- New metadata rows 4-10 ---
$ws.Range("A4").Value = "dct:creator"
$ws.Range("B4").Value = "Mario Scrocca (Cefriel)"

$ws.Range("A5").Value = "dct:publisher"
$ws.Range("B5").Value = "TANGENT WP2"

$ws.Range("A6").Value = "owl:versionInfo"
$ws.Range("B6").Value = "1.0.0"

$ws.Range("A7").Value = "owl:versionIRI"
$ws.Range("B7").Formula = "=_xlfn.CONCAT(B1,""/"",B6)"

$ws.Range("A8").Value = "owl:priorVersion"

$ws.Range("A9").Value = "dct:license"

$ws.Range("A10").Value = "http://purl.org/ontology/bibo/status"
$ws.Range("B10").Value = "Published Controlled Vocabulary"
# Row 10 has no formatting in column C in the target layout.
$ws.Range("C10").Clear()

# --- Apply "Hyperlink" look to B8 (left empty) and B9, then make B9 a real hyperlink ---
$b8 = $ws.Range("B8")
$b8.ClearFormats()
$b8.Value = "temp-placeholder"
$tempLink = $ws.Hyperlinks.Add($b8, "https://example.com/temp-placeholder")
$tempLink.Delete()
$b8.Value = $null

$b9 = $ws.Range("B9")
$b9.ClearFormats()
$b9.Value = "https://creativecommons.org/licenses/by/4.0/"
$ws.Hyperlinks.Add($b9, "https://creativecommons.org/licenses/by/4.0/") | Out-Null

# Match the reported selection state of the saved worksheet.
$ws.Range("A1:XFD10").Select() | Out-Null

$wb.Save()
